# MeetTeam.pptx — slide 1: swap the Title/Subtitle copy.
#   Title   : "Meet Your Team 👋"  ->  "App Project 2024"
#   Subtitle: "App Project 2024"   ->  "Meet Your Team 👋" (two runs: "Meet " / "Your Team 👋")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title (shape 1: "Title 1") ------------------------------------------
# Fully replace the title text. Deleting first + inserting fresh text collapses
# the paragraph to a single run (taking on the first run's plain rPr) and drops
# the now-unneeded endParaRPr, matching the target markup exactly.
$titleShape = $s.Shapes.Item("Title 1")
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Delete()
$null = $titleRange.InsertAfter("App Project 2024")

# --- Subtitle (shape 2: "Subtitle 2") -------------------------------------
# Replace with the combined text first (keeps the existing run's rich
# formatting: solidFill 333333, white highlight, Inter latin typeface, and
# the paragraph's endParaRPr), then re-touch the Font of the trailing
# "Your Team 👋" span so it becomes its own run while keeping identical
# formatting to "Meet ".
$subShape = $s.Shapes.Item("Subtitle 2")
$subRange = $subShape.TextFrame.TextRange
$subRange.Text = "Meet Your Team 👋"

$tail = $subRange.Characters(6, 11)
$tail.Font.Color.RGB = $tail.Font.Color.RGB
$null = $tail
